$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)
$s.Delete()
